# "Penalty Reward System" (unfinished) — the forecast got re-pulled a week
# later than before, so every week's start date shifts forward by one week
# (the series now runs 2025-01-12 .. 2025-04-27 instead of
# 2025-01-05 .. 2025-04-20) and MyForecast for each of those weeks was
# overwritten by the new (buggy/unfinished) penalty-reward values.
#
# NOTE: Week_Start_Date / the Summary "Value" column are plain text
# (inlineStr) in the source file, not real dates — several look like dates
# or bare numbers ("218", "2025-03-23", ...). Setting .Value directly on a
# string that LOOKS like a date/number makes Excel auto-coerce it to a
# date serial / number. Forcing NumberFormat to Text ("@") before the
# write keeps it literal text, and resetting Style back to "Normal"
# afterwards drops the leftover number-format so the cell ends up with no
# style override at all (matching the original formatting).

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Forecast Comparison")
$ws2 = $wb.Worksheets.Item("Summary")

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# --- Forecast Comparison (rows 2..17): Week_Start_Date shifts one week
#     forward; MyForecast collapses to the new 13/14 penalty-reward values.
$weekStarts = @(
    "2025-01-12", "2025-01-19", "2025-01-26", "2025-02-02",
    "2025-02-09", "2025-02-16", "2025-02-23", "2025-03-02",
    "2025-03-09", "2025-03-16", "2025-03-23", "2025-03-30",
    "2025-04-06", "2025-04-13", "2025-04-20", "2025-04-27"
)
$myForecasts = @(13, 13, 13, 14, 14, 14, 14, 14, 14, 14, 14, 14, 14, 14, 14, 14)

for ($i = 0; $i -lt $weekStarts.Count; $i++) {
    $row = $i + 2
    Set-TextValue $ws1.Range("B$row") $weekStarts[$i]
    $ws1.Range("D$row").Value = $myForecasts[$i]
}

# --- Summary sheet: recomputed (unfinished) rollup figures ---
Set-TextValue $ws2.Range("B2")  "2023-12-31 to 2025-01-05"
Set-TextValue $ws2.Range("B8")  "664 units"
Set-TextValue $ws2.Range("B9")  "218"
Set-TextValue $ws2.Range("B10") "108"
Set-TextValue $ws2.Range("B11") "53"
Set-TextValue $ws2.Range("B12") "14"
Set-TextValue $ws2.Range("B13") "2025-03-23"
Set-TextValue $ws2.Range("B14") "13"
Set-TextValue $ws2.Range("B15") "2025-01-12"
